$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.8703656764443379, 1.0320466954387313, 0.5465517152907997, 2.2207479399061842)"
$ws.Range("C2").Value = "NIG(1.5800297202085303, 1.1741681119539105, 3.6198597448984797, 6.01648817107456)"
$ws.Range("D2").Value = "JSU(-0.8192558337956347, 1.2788509834542616, 1.4195791808988019, 3.5271215308500636)"
$ws.Range("E2").Value = "EXN(1.7893037226222037, 3.408661499577847, 3.900033340598612)"
